$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current row 297 (shifts old 297-306 down to 299-308)
$ws.Range("A297:A298").EntireRow.Insert()

# Copy the date-formatted style from a known formatted date cell (D299, which is the
# old D297 pushed down) onto the two new date cells so they keep the same number format.
$ws.Range("D299").Copy()
$ws.Range("D297:D298").PasteSpecial(-4122)  # xlPasteFormats

# New row 297 data
$ws.Cells.Item(297, 1).Value = 2
$ws.Cells.Item(297, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(297, 3).Value = "Coquimbo"
$ws.Cells.Item(297, 4).Value = 45267
$ws.Cells.Item(297, 5).Value = 4
$ws.Cells.Item(297, 6).Value = 100112031
$ws.Cells.Item(297, 7).Value = "Poroto verde"
$ws.Cells.Item(297, 8).Value = "Magnum"
$ws.Cells.Item(297, 9).Value = "Primera"
$ws.Cells.Item(297, 10).Value = 800
$ws.Cells.Item(297, 11).Value = 27000
$ws.Cells.Item(297, 12).Value = 30000
$ws.Cells.Item(297, 13).Value = 28500
$ws.Cells.Item(297, 14).Value = "$/caja 25 kilos"
$ws.Cells.Item(297, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(297, 16).Value = 1140
$ws.Cells.Item(297, 17).Value = 25
$ws.Cells.Item(297, 18).Value = "Hortaliza"

# New row 298 data
$ws.Cells.Item(298, 1).Value = 2
$ws.Cells.Item(298, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(298, 3).Value = "Coquimbo"
$ws.Cells.Item(298, 4).Value = 45267
$ws.Cells.Item(298, 5).Value = 4
$ws.Cells.Item(298, 6).Value = 100112031
$ws.Cells.Item(298, 7).Value = "Poroto verde"
$ws.Cells.Item(298, 8).Value = "Sin especificar"
$ws.Cells.Item(298, 9).Value = "Primera"
$ws.Cells.Item(298, 10).Value = 800
$ws.Cells.Item(298, 11).Value = 30000
$ws.Cells.Item(298, 12).Value = 32000
$ws.Cells.Item(298, 13).Value = 31000
$ws.Cells.Item(298, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(298, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(298, 16).Value = 1240
$ws.Cells.Item(298, 17).Value = 25
$ws.Cells.Item(298, 18).Value = "Hortaliza"
